$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "value" -> "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# Propagate the date-column style (A2) down to the new rows so they
# all share the same cell style (border/font/number format) instead of
# Excel fabricating a brand-new style entry per cell.
$ws.Range("A2").Copy($ws.Range("A3:A22"))

# Data rows: (row, dateSerial, value-or-$null)
$data = @(
    @(2, 38717, $null),
    @(3, 39082, -0.1287148287979267),
    @(4, 39447, 0.08844991283951664),
    @(5, 39813, 0.9225722794137248),
    @(6, 40178, -0.7475385776494314),
    @(7, 40543, 1.094490700431927),
    @(8, 40908, -0.8742167833903691),
    @(9, 41274, -0.5751029748885195),
    @(10, 41639, -0.1263097576649996),
    @(11, 42004, 0.4976690624053814),
    @(12, 42369, -0.3267144271395628),
    @(13, 42735, 0.149524011641855),
    @(14, 43100, 0.1983963998054783),
    @(15, 43465, 0.6497679376401333),
    @(16, 43830, -0.5803176690338252),
    @(17, 44196, 0.2929419500579789),
    @(18, 44561, -2.411738983600742),
    @(19, 44926, -1.022826494952023),
    @(20, 45291, -0.5858537819409149),
    @(21, 45657, 0.01773204329378331),
    @(22, 46022, $null)
)

foreach ($item in $data) {
    $row = $item[0]
    $dateSerial = $item[1]
    $value = $item[2]

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $dateSerial

    $bCell = $ws.Cells.Item($row, 2)
    if ($null -ne $value) {
        $bCell.Value = $value
    } else {
        $bCell.Value = $null
    }
}
